# Apply the Alvearie FHIR IG "total-dependent-count" StructureDefinition
# refresh (5.0.0 -> 6.0.0) to the workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Metadata" ----
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Date refresh
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now populated
$meta.Range("B9").Value = "Alvearie Team"

# Old duplicate "Contact" row (row 10) becomes the real "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Remove the leftover duplicate "Contact" / "No display for ContactDetail" row
$meta.Rows.Item(11).Delete()

# ---- Sheet "Elements" ----
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the top-level "Extension" element; give it the real short/definition text
$elements.Range("K2").Value = "Total Dependent Count"
$elements.Range("L2").Value = "Number of dependents covered by the group health medical plan"
